$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1 (15:22 -> 15:52)
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 15:52"

# Update the Bizkaia/Vizcaya row (row 6) statistics
$ws.Range("B6").Value = 6538
$ws.Range("C6").Value = 800
$ws.Range("D6").Value = 5460
$ws.Range("E6").Value = 278
